$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.649.34"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "1.584.45"
$ws.Range("E3").Value = "  -3.15%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.57"
$ws.Range("E5").Value = "  -2.44%  "
$ws.Range("E6").Value = "  -3.25%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  -4.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.253"
$ws.Range("E9").Value = "  -1.28%  "
$ws.Range("E10").Value = "  -3.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0869"
$ws.Range("E11").Value = "  -1.63%  "
$ws.Range("D12").Value = "1.809.92"
$ws.Range("E12").Value = "  -3.12%  "
$ws.Range("D13").Value = "1.586.86"
$ws.Range("E13").Value = "  -3.12%  "
$ws.Range("E14").Value = "  -3.92%  "
$ws.Range("E15").Value = "  -5.71%  "
$ws.Range("D16").Value = "27.647.99"
$ws.Range("E16").Value = "  -0.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.12"
$ws.Range("E17").Value = "  -3.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "219.21"
$ws.Range("E18").Value = "  -4.17%  "
$ws.Range("E19").Value = "  -3.51%  "
$ws.Range("E20").Value = "  -5.04%  "
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("E22").Value = "  -5.01%  "
$ws.Range("E23").Value = "  -5.33%  "
$ws.Range("E24").Value = "  -5.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.56"
$ws.Range("E25").Value = "  -1.07%  "
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("E27").Value = "  -2.49%  "
$ws.Range("E28").Value = "  -2.90%  "
$ws.Range("E29").Value = "  -4.08%  "
$ws.Range("E30").Value = "  -2.82%  "
$ws.Range("E31").Value = "  -3.27%  "
$ws.Range("E32").Value = "  -5.38%  "
$ws.Range("D33").Value = "1.382.38"
$ws.Range("E33").Value = "  -1.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.96"
$ws.Range("E34").Value = "  -4.81%  "
$ws.Range("E35").Value = "  -5.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.963"
$ws.Range("E36").Value = "  -4.77%  "
$ws.Range("E37").Value = "  -1.08%  "
$ws.Range("E38").Value = "  -2.83%  "
$ws.Range("E39").Value = "  -3.17%  "
$ws.Range("E40").Value = "  -3.79%  "
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.979"
$ws.Range("E42").Value = "  -2.48%  "
$ws.Range("E43").Value = "  -3.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.77"
$ws.Range("E44").Value = "  -3.55%  "
$ws.Range("E45").Value = "  +1.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.24"
$ws.Range("E46").Value = "  -4.10%  "
$ws.Range("D47").Value = "1.721.74"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.15"
$ws.Range("E48").Value = "  -0.77%  "
$ws.Range("D49").Value = "0.0₆0101"
$ws.Range("E49").Value = "  -1.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0975"
$ws.Range("E50").Value = "  -5.21%  "
$ws.Range("E51").Value = "  -1.15%  "
